# Update the crawl timestamp column (O) for all data rows (2-549)
# from "2022-12-23 06:52:13" to "2022-12-23 12:58:44",
# and update the product aria-label text in M333 to reflect the
# "Online kein Bestand" (no online stock) note for the Betty Bossi
# Mungbohnen Sprossen product.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2022-12-23 06:52:13"
$newTimestamp = "2022-12-23 12:58:44"

$lastRow = 549
$firstRow = 2

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 15)  # column O = 15
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}

# Row 333, column M (13): update product aria label text
$ws.Cells.Item(333, 13).Value = "Betty Bossi Mungbohnen Sprossen - Online kein Bestand 2.30 Schweizer Franken"
